$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the "_old" / "_new" header labels to "_FV2404" / "_FV2410".
#    (Row 1 holds the column headers; "diff" in K1 is left untouched.)
# ---------------------------------------------------------------------------
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U55 into a real Excel table ("Table1") without letting the
#    table-creation step invent a header dxf: Excel only emits
#    headerRowDxfId / a styles.xml <dxf> when the header range carries direct
#    formatting at the moment the ListObject is created, so we stash the
#    header's current formatting, blank it out, add the table, then restore
#    the formatting via copy/paste (which reuses the existing cellXfs entry
#    instead of minting a differential-format record).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$backupRange = $ws.Range("A1000:U1000")

$headerRange.Copy($backupRange)
$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$backupRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$backupRange.Clear()

$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split above row 2, pane frozen).
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
